$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = 45224
$ws.Cells.Item(2, 13).Value = 80
$ws.Cells.Item(2, 14).Value = 20000
$ws.Cells.Item(2, 15).Value = 20000
$ws.Cells.Item(2, 16).Value = 20000
$ws.Cells.Item(2, 19).Value = 2000

$ws.Cells.Item(3, 4).Value = 44447
$ws.Cells.Item(3, 13).Value = 60
$ws.Cells.Item(3, 14).Value = 21000
$ws.Cells.Item(3, 15).Value = 22000
$ws.Cells.Item(3, 16).Value = 21500
$ws.Cells.Item(3, 19).Value = 2150

$ws.Cells.Item(4, 4).Value = 45203
$ws.Cells.Item(4, 13).Value = 30
$ws.Cells.Item(4, 14).Value = 21000
$ws.Cells.Item(4, 15).Value = 21000
$ws.Cells.Item(4, 16).Value = 21000
$ws.Cells.Item(4, 19).Value = 2100

$ws.Cells.Item(5, 4).Value = 45205
$ws.Cells.Item(5, 14).Value = 22000
$ws.Cells.Item(5, 15).Value = 22000
$ws.Cells.Item(5, 16).Value = 22000
$ws.Cells.Item(5, 19).Value = 2200

$ws.Cells.Item(6, 4).Value = 44848
$ws.Cells.Item(6, 12).Value = "Especial"
$ws.Cells.Item(6, 13).Value = 60
$ws.Cells.Item(6, 14).Value = 24000
$ws.Cells.Item(6, 15).Value = 25000
$ws.Cells.Item(6, 16).Value = 24500
$ws.Cells.Item(6, 19).Value = 2450

$ws.Cells.Item(7, 12).Value = "Primera"
$ws.Cells.Item(7, 13).Value = 120
$ws.Cells.Item(7, 14).Value = 21000
$ws.Cells.Item(7, 15).Value = 22000
$ws.Cells.Item(7, 16).Value = 21500
$ws.Cells.Item(7, 19).Value = 2150

$ws.Cells.Item(8, 4).Value = 44448
$ws.Cells.Item(8, 13).Value = 60

$ws.Cells.Item(9, 4).Value = 45180
$ws.Cells.Item(9, 12).Value = "Primera"
$ws.Cells.Item(9, 13).Value = 40
$ws.Cells.Item(9, 14).Value = 22000
$ws.Cells.Item(9, 15).Value = 22000
$ws.Cells.Item(9, 16).Value = 22000
$ws.Cells.Item(9, 19).Value = 2200

$ws.Cells.Item(10, 12).Value = "Especial"
$ws.Cells.Item(10, 13).Value = 60
$ws.Cells.Item(10, 14).Value = 31000
$ws.Cells.Item(10, 15).Value = 32000
$ws.Cells.Item(10, 16).Value = 31500
$ws.Cells.Item(10, 19).Value = 3150

$ws.Cells.Item(11, 4).Value = 44460
$ws.Cells.Item(11, 13).Value = 30
$ws.Cells.Item(11, 14).Value = 30000
$ws.Cells.Item(11, 15).Value = 30000
$ws.Cells.Item(11, 16).Value = 30000
$ws.Cells.Item(11, 19).Value = 3000

$ws.Cells.Item(12, 4).Value = 44874
$ws.Cells.Item(12, 12).Value = "Especial"
$ws.Cells.Item(12, 13).Value = 30
$ws.Cells.Item(12, 14).Value = 25000
$ws.Cells.Item(12, 15).Value = 25000
$ws.Cells.Item(12, 16).Value = 25000
$ws.Cells.Item(12, 19).Value = 2500

$ws.Cells.Item(13, 4).Value = 44874
$ws.Cells.Item(13, 13).Value = 80
$ws.Cells.Item(13, 14).Value = 23000
$ws.Cells.Item(13, 15).Value = 24000
$ws.Cells.Item(13, 16).Value = 23500
$ws.Cells.Item(13, 19).Value = 2350

$ws.Cells.Item(14, 4).Value = 45191
$ws.Cells.Item(14, 13).Value = 30
$ws.Cells.Item(14, 14).Value = 21000
$ws.Cells.Item(14, 15).Value = 21000
$ws.Cells.Item(14, 16).Value = 21000
$ws.Cells.Item(14, 19).Value = 2100

$ws.Cells.Item(15, 4).Value = 45189
$ws.Cells.Item(15, 13).Value = 50
$ws.Cells.Item(15, 14).Value = 22000
$ws.Cells.Item(15, 15).Value = 22000
$ws.Cells.Item(15, 16).Value = 22000
$ws.Cells.Item(15, 19).Value = 2200

$ws.Cells.Item(16, 4).Value = 44839
$ws.Cells.Item(16, 13).Value = 120
$ws.Cells.Item(16, 14).Value = 25000
$ws.Cells.Item(16, 15).Value = 26000
$ws.Cells.Item(16, 16).Value = 25500
$ws.Cells.Item(16, 19).Value = 2550

$ws.Cells.Item(17, 4).Value = 45225
$ws.Cells.Item(17, 13).Value = 80
$ws.Cells.Item(17, 14).Value = 21000
$ws.Cells.Item(17, 15).Value = 21000
$ws.Cells.Item(17, 16).Value = 21000
$ws.Cells.Item(17, 19).Value = 2100

$ws.Cells.Item(18, 4).Value = 45212
$ws.Cells.Item(18, 12).Value = "Primera"
$ws.Cells.Item(18, 13).Value = 80
$ws.Cells.Item(18, 14).Value = 22000
$ws.Cells.Item(18, 15).Value = 22000
$ws.Cells.Item(18, 16).Value = 22000
$ws.Cells.Item(18, 19).Value = 2200

$ws.Cells.Item(19, 4).Value = 44868
$ws.Cells.Item(19, 12).Value = "Especial"
$ws.Cells.Item(19, 13).Value = 60
$ws.Cells.Item(19, 14).Value = 26000
$ws.Cells.Item(19, 15).Value = 26000
$ws.Cells.Item(19, 16).Value = 26000
$ws.Cells.Item(19, 19).Value = 2600

$ws.Cells.Item(20, 4).Value = 45194
$ws.Cells.Item(20, 12).Value = "Primera"
$ws.Cells.Item(20, 13).Value = 80
$ws.Cells.Item(20, 14).Value = 22000
$ws.Cells.Item(20, 15).Value = 22000
$ws.Cells.Item(20, 16).Value = 22000
$ws.Cells.Item(20, 19).Value = 2200

$ws.Cells.Item(21, 4).Value = 44487
$ws.Cells.Item(21, 14).Value = 23000
$ws.Cells.Item(21, 15).Value = 24000
$ws.Cells.Item(21, 16).Value = 23500
$ws.Cells.Item(21, 19).Value = 2350

$ws.Cells.Item(23, 4).Value = 45216
$ws.Cells.Item(23, 13).Value = 60

$ws.Cells.Item(24, 4).Value = 45196
$ws.Cells.Item(24, 14).Value = 23000
$ws.Cells.Item(24, 15).Value = 23000
$ws.Cells.Item(24, 16).Value = 23000
$ws.Cells.Item(24, 19).Value = 2300

$ws.Cells.Item(25, 4).Value = 45230
$ws.Cells.Item(25, 14).Value = 21000
$ws.Cells.Item(25, 15).Value = 21000
$ws.Cells.Item(25, 16).Value = 21000
$ws.Cells.Item(25, 19).Value = 2100

$ws.Cells.Item(26, 4).Value = 44452
$ws.Cells.Item(26, 13).Value = 60
$ws.Cells.Item(26, 14).Value = 21000
$ws.Cells.Item(26, 16).Value = 21500
$ws.Cells.Item(26, 19).Value = 2150

$ws.Cells.Item(27, 4).Value = 45217
$ws.Cells.Item(27, 14).Value = 21000
$ws.Cells.Item(27, 15).Value = 21000
$ws.Cells.Item(27, 16).Value = 21000
$ws.Cells.Item(27, 19).Value = 2100

$ws.Cells.Item(28, 4).Value = 45176
$ws.Cells.Item(28, 12).Value = "Primera"
$ws.Cells.Item(28, 14).Value = 22000
$ws.Cells.Item(28, 15).Value = 22000
$ws.Cells.Item(28, 16).Value = 22000
$ws.Cells.Item(28, 19).Value = 2200

$ws.Cells.Item(29, 4).Value = 44841
$ws.Cells.Item(29, 13).Value = 60

$ws.Cells.Item(30, 4).Value = 45209
$ws.Cells.Item(30, 13).Value = 50
$ws.Cells.Item(30, 14).Value = 22000
$ws.Cells.Item(30, 16).Value = 22000
$ws.Cells.Item(30, 19).Value = 2200

$ws.Cells.Item(31, 4).Value = 45219
$ws.Cells.Item(31, 14).Value = 20000
$ws.Cells.Item(31, 15).Value = 20000
$ws.Cells.Item(31, 16).Value = 20000
$ws.Cells.Item(31, 19).Value = 2000

$ws.Cells.Item(32, 4).Value = 44461
$ws.Cells.Item(32, 12).Value = "Especial"
$ws.Cells.Item(32, 14).Value = 31000
$ws.Cells.Item(32, 15).Value = 32000
$ws.Cells.Item(32, 16).Value = 31500
$ws.Cells.Item(32, 19).Value = 3150

$ws.Cells.Item(33, 4).Value = 44461
$ws.Cells.Item(33, 13).Value = 30
$ws.Cells.Item(33, 14).Value = 30000
$ws.Cells.Item(33, 15).Value = 30000
$ws.Cells.Item(33, 16).Value = 30000
$ws.Cells.Item(33, 19).Value = 3000

$ws.Cells.Item(34, 4).Value = 45173
$ws.Cells.Item(34, 13).Value = 50
$ws.Cells.Item(34, 14).Value = 22000
$ws.Cells.Item(34, 15).Value = 22000
$ws.Cells.Item(34, 16).Value = 22000
$ws.Cells.Item(34, 19).Value = 2200
